$wb = $excel.ActiveWorkbook

# "Generate Report for Archive": the status text "Ready for handoff" is now
# "In Translation" everywhere it appears, and the Status column is
# re-autofit (narrower) to match the new, shorter text.
$newStatus = "In Translation"
$newColumnWidth = 12.5

# Overview sheet: columns E (zh-cn) and F (de-de) show each locale's status
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $newColumnWidth
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth

# zh-cn sheet: column C is the Status column
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth

# de-de sheet: column C is the Status column
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth
